# "Generate Report for Handback"
#
# This localization-status report is regenerated on every handback run:
# the per-language Status + the "Latest Handback DateTime" / "Error Detail"
# columns get refreshed for each language sheet, and the Overview sheet's
# per-language status cells follow suit.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Sheets.Item("Overview")
$ws_zhcn     = $wb.Sheets.Item("zh-cn")
$ws_dede     = $wb.Sheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: both language status columns move to "handed back" ---
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus

# --- zh-cn sheet: handback timestamp refreshed, handback is in sync (no error) ---
$ws_zhcn.Range("K2").Value = "2016-08-14 17:06:14"
$ws_zhcn.Range("P2").Value = ""

# --- de-de sheet: handback timestamp refreshed, handback is in sync (no error) ---
$ws_dede.Range("K2").Value = "2016-08-14 17:06:24"
$ws_dede.Range("P2").Value = ""

# --- Column widths re-fit to the new cell content ---
$ws_overview.Columns.Item(5).AutoFit()
$ws_overview.Columns.Item(6).AutoFit()

$ws_zhcn.Columns.Item(3).AutoFit()
$ws_zhcn.Columns.Item(16).AutoFit()

$ws_dede.Columns.Item(3).AutoFit()
$ws_dede.Columns.Item(16).AutoFit()
